$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -12.679
    10 = -12.805
    12 = -11.066
    18 = -11.409
    37 = -13.049
    55 = -13.637
    68 = -11.022
    77 = -12.824
    78 = -12.798
    81 = -13.094
    82 = -11.885
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
